$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SoCtMbCtbDP")
$ws.Range("B1:B24").Value = 0.8
